$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Add the two new ETF rows (Gold & Silver) at the bottom of the table
# ---------------------------------------------------------------

# Row 32 - ICICI Prudential Gold ETF
$ws.Range("A32").Value = 33
$ws.Range("B32").Value = "ICICI Prudential Gold ETF"
$ws.Range("C32").Value = "GOLDIETF"
$ws.Range("D32").Value = 201
$ws.Range("E32").Value = 123.86

# Row 33 - ICICI Prudential Silver ETF
$ws.Range("A33").Value = 34
$ws.Range("B33").Value = "ICICI Prudential Silver ETF"
$ws.Range("C33").Value = "SILVERIETF"
$ws.Range("D33").Value = 108
$ws.Range("E33").Value = 229.98

# ---------------------------------------------------------------
# Formatting to match the rest of the table (thin black borders all
# around, "Aptos Narrow" font + text format on the asset-name column)
# ---------------------------------------------------------------

$dataRange = $ws.Range("A32:E33")
$dataRange.Borders.ColorIndex = 1
$dataRange.Borders.LineStyle = 1
$dataRange.Borders.Weight = 2

$assetRange = $ws.Range("B32:B33")
$assetRange.Font.Name = "Aptos Narrow"
$assetRange.NumberFormat = "@"

# ---------------------------------------------------------------
# Sheet view bookkeeping (matches final selection/scroll position)
# ---------------------------------------------------------------

$ws.Range("I31").Select()
$ws.Application.ActiveWindow.ScrollRow = 13

Write-Host "done"
